$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update Meta/Venda values
$ws.Range("A5").Value = 6800
$ws.Range("B5").Value = 9000

# Row 32: convert inline-string text values to real numbers
$ws.Range("A32").Value = 35600
$ws.Range("B32").Value = 64641
